$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_7_6_0",0.7910188390328098,0.9091011449515201,0.6254767037754212,0.8229452383691889,0.231280580163002,0.07942944765090942,0.1613206267356873,0.1179664954543114),
    @("model_7_6_1",0.7981643045777636,0.9060960606875952,0.632906261595537,0.8231190293777405,0.2233726531267166,0.08205536007881165,0.1581204533576965,0.1178507208824158),
    @("model_7_6_2",0.804942623137531,0.9030424379709416,0.6389978019014303,0.8228522691928117,0.2158710658550262,0.08472368121147156,0.1554965972900391,0.1180284470319748),
    @("model_7_6_3",0.8113253129270077,0.9001387259768857,0.6428252023006018,0.8220002995581144,0.2088072597980499,0.08726099878549576,0.1538480073213577,0.1185961067676544),
    @("model_7_6_4",0.8174161782298359,0.897949950368048,0.6447315863900929,0.8210605209014048,0.202066496014595,0.08917360752820969,0.1530268490314484,0.1192222461104393),
    @("model_7_6_5",0.8232343318184314,0.8969348752835573,0.6437418322543139,0.8200546170853158,0.1956274956464767,0.09006060659885406,0.1534532010555267,0.1198924407362938),
    @("model_7_6_6",0.8286929913456262,0.8966629888947315,0.6392254365431521,0.8184919409160585,0.1895864009857178,0.09029818326234818,0.1553985625505447,0.1209336146712303),
    @("model_7_6_7",0.8336196714084425,0.897269487860074,0.6270836326668295,0.8152191070773965,0.1841340065002441,0.08976820856332779,0.160628467798233,0.1231142058968544),
    @("model_7_6_8",0.8376790581714322,0.8974230997618156,0.6050047227844371,0.8086088410123528,0.179641455411911,0.08963398635387421,0.1701386570930481,0.127518430352211),
    @("model_7_6_9",0.8413862073649325,0.8970569662045559,0.5816024496113248,0.8012347208523009,0.1755387336015701,0.08995392173528671,0.1802188456058502,0.1324316114187241),
    @("model_7_6_10",0.8445750166408454,0.8969238451907887,0.5518016466609983,0.7920759897466756,0.1720096617937088,0.09007023274898529,0.193055123090744,0.1385337859392166),
    @("model_7_6_11",0.8469205605554602,0.8974631122924752,0.5096265639152572,0.7796194840225399,0.169413834810257,0.08959902822971344,0.21122145652771,0.1468332260847092),
    @("model_7_6_12",0.8497007450963017,0.8980076218657798,0.4784574437634944,0.7705150278482567,0.1663369834423065,0.08912321925163269,0.2246471047401428,0.1528992503881454),
    @("model_7_6_13",0.8526443477904422,0.89931909943212,0.4511911766722269,0.7631304698514239,0.1630792766809464,0.08797720819711685,0.2363916635513306,0.1578193753957748),
    @("model_7_6_14",0.8558885954418365,0.9004775348277246,0.4323741572685847,0.7582102774681413,0.1594888716936111,0.08696494996547699,0.2444968223571777,0.1610975563526154),
    @("model_7_6_15",0.8587544721217211,0.9020175742707286,0.4100169310387078,0.7524776572269817,0.1563171595335007,0.08561922609806061,0.2541269063949585,0.1649170517921448),
    @("model_7_6_16",0.8619603913881422,0.9032783754109351,0.3964660673903179,0.7492306994668286,0.1527691632509232,0.08451750874519348,0.2599637508392334,0.1670803874731064),
    @("model_7_6_17",0.8648180123369733,0.9039126083215101,0.3810093083093736,0.744968629815292,0.1496066004037857,0.0839633047580719,0.2666215300559998,0.1699200868606567),
    @("model_7_6_18",0.8674981619950888,0.9042216175821491,0.3656652960804211,0.7405151142652997,0.1466404944658279,0.08369328826665878,0.2732307314872742,0.1728873252868652),
    @("model_7_6_19",0.8699387760751471,0.9043637920244048,0.3485859248851034,0.7354177920166995,0.1439394652843475,0.08356904983520508,0.2805874347686768,0.1762835085391998),
    @("model_7_6_20",0.8723083844495988,0.9048619677884445,0.3313602801925823,0.7305231386749073,0.1413169950246811,0.08313372731208801,0.2880071401596069,0.1795447021722794),
    @("model_7_6_21",0.8744311029510601,0.9050276346345814,0.313396565947255,0.7251729275811785,0.1389677673578262,0.08298896998167038,0.2957447171211243,0.1831094026565552),
    @("model_7_6_22",0.876350375712641,0.9041097986944752,0.2964842342606882,0.7193906941539605,0.1368436962366104,0.08379099518060684,0.303029477596283,0.1869619190692902),
    @("model_7_6_23",0.8780612108997591,0.9034937978570214,0.2771240952094969,0.7130729002034555,0.1349503099918365,0.0843292772769928,0.3113686144351959,0.1911712884902954),
    @("model_7_6_24",0.879625857979323,0.9026051856181353,0.2574701492143713,0.7064766014869475,0.1332186907529831,0.08510575443506241,0.3198342323303223,0.1955662220716476)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item(2 + $i, 1 + $j).Value = $rowData[$j]
    }
}
